# Update report for Palagano to 23 August 2021 (aggiornamento al 23 agosto 2021)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data ends at row 343 (date serial 44417, i.e. 2021-08-09).
# Append rows 344..357 (dates 44418..44431, i.e. 2021-08-10..2021-08-23).

$lastRow = 343
$firstNewRow = 344
$lastNewRow = 357

# Copy the date-column formatting (style) from the last existing data row
# down through the new rows, matching the border/bold/centered date style.
$ws.Range("A$lastRow").Copy($ws.Range("A${firstNewRow}:A${lastNewRow}"))

# Values for B (nuovi pos.), C (somma mobile 7gg.), D (somma mobile 7gg. per 100mila ab.)
$values = @{
    344 = @(0, 0, 0)
    345 = @(0, 0, 0)
    346 = @(0, 0, 0)
    347 = @(0, 0, 0)
    348 = @(0, 0, 0)
    349 = @(0, 0, 0)
    350 = @(0, 0, 0)
    351 = @(0, 0, 0)
    352 = @(0, 0, 0)
    353 = @(0, 0, 0)
    354 = @(0, 0, 0)
    355 = @(0, 0, 0)
    356 = @(2, 2, 96.15384615384616)
    357 = @(0, 2, 96.15384615384616)
}

$serial = 44418
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $serial
    $v = $values[$r]
    $ws.Cells.Item($r, 2).Value = $v[0]
    $ws.Cells.Item($r, 3).Value = $v[1]
    $ws.Cells.Item($r, 4).Value = $v[2]
    $serial = $serial + 1
}
